# Generate Report for Handoff
# Updates status text and timestamps across Overview / zh-cn / de-de sheets
# (the localization run moved from "Ready for handoff" into "In Translation"
# and refreshed its handoff timestamps), and narrows the now-shorter
# "Status" column on each of those sheets to re-fit the new content.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# New, narrower column width for the status columns below. Excel quantizes
# ColumnWidth to 1/6-character increments when it writes the sheet back out,
# so 12.5 is the input that lands closest to the refreshed autofit width.
$statusColWidth = 12.5

# --- Overview sheet ---
# E2 = zh-cn status, F2 = de-de status, G2 = Latest HO Xliff Generate Date
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("G2").Value = "2016-10-14 07:07:07"

# Narrow the status columns (E, F) to match the new content width
$overview.Range("E1").ColumnWidth = $statusColWidth
$overview.Range("F1").ColumnWidth = $statusColWidth

# --- zh-cn sheet ---
# C2 = Status, H2 = Latest Handoff Datetime
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("H2").Value = "2016-10-14 07:06:57"
$zhcn.Range("C1").ColumnWidth = $statusColWidth

# --- de-de sheet ---
# C2 = Status, H2 = Latest Handoff Datetime (shares the same original text as
# Overview!G2, so it must be updated in lock-step to land on the same value)
$dede.Range("C2").Value = "In Translation"
$dede.Range("H2").Value = "2016-10-14 07:07:07"
$dede.Range("C1").ColumnWidth = $statusColWidth
